# "Edit report template, fix submit thesis"
# The approved-thesis export template gains a "Giảng viên ra đề"
# (lecturer who set the topic) column, and the sheet/tab name is updated
# to reflect that this is the list of APPROVED theses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: "DANH SÁCH ĐỀ TÀI" -> "DANH SÁCH ĐỀ TÀI ĐÃ ĐƯỢC DUYỆT"
$ws.Name = "DANH SÁCH ĐỀ TÀI ĐÃ ĐƯỢC DUYỆT"

# Insert a new column E ("Giảng viên ra đề" / "{{Items.LecturerName}}"),
# shifting the old E..I columns (credits..specialization) right to F..J.
$ws.Columns("E").Insert()

# New column E header (row2, matches the other header cells' style) and
# template placeholder (row3, matches the other body cells' style).
$ws.Range("F2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = "Giảng viên ra đề"

$ws.Range("H3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "{{Items.LecturerName}}"

$excel.CutCopyMode = $false

# Column widths: new column E (Lecturer) and the shifted-right last column J
# (Specialization) both got new, wider custom widths.
$ws.Columns("E").ColumnWidth = 29.5
$ws.Columns("J").ColumnWidth = 28.6

# Column G (old F, "Năm thực hiện") keeps its original width but is now its
# own <col> record rather than merged with column F, since column F's width
# changed independently when column E's width was set above - re-touch it
# so it is written out as a distinct run.
$ws.Columns("G").ColumnWidth = $ws.Columns("G").ColumnWidth

# Selection moved in the saved view.
$ws.Range("C14").Select()
